$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly pair of rows (date 2023-07-17 / serial 45124) is inserted right
# before the existing row 1040, shifting all the data that follows down by two
# rows (old row 1170 -> new row 1172; sheet dimension A1:R1170 -> A1:R1172).
$ws.Rows("1040:1041").Insert()

# Row 1040: "Primera" quality for the new date.
$ws.Range("A1040").Value = 8
$ws.Range("B1040").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1040").Value = "Coquimbo"
$ws.Range("D1040").Value = 45124
$ws.Range("E1040").Value = 4
$ws.Range("F1040").Value = 100112023
$ws.Range("G1040").Value = "Brócoli"
$ws.Range("H1040").Value = "Sin especificar"
$ws.Range("I1040").Value = "Primera"
$ws.Range("J1040").Value = 2000
$ws.Range("K1040").Value = 700
$ws.Range("L1040").Value = 800
$ws.Range("M1040").Value = 750
$ws.Range("N1040").Value = "$/unidad"
$ws.Range("O1040").Value = "Provincia del Elquí"
$ws.Range("P1040").Value = 750
$ws.Range("Q1040").Value = 1
$ws.Range("R1040").Value = "Hortaliza"

# Row 1041: "Segunda" quality for the new date.
$ws.Range("A1041").Value = 8
$ws.Range("B1041").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1041").Value = "Coquimbo"
$ws.Range("D1041").Value = 45124
$ws.Range("E1041").Value = 4
$ws.Range("F1041").Value = 100112023
$ws.Range("G1041").Value = "Brócoli"
$ws.Range("H1041").Value = "Sin especificar"
$ws.Range("I1041").Value = "Segunda"
$ws.Range("J1041").Value = 1200
$ws.Range("K1041").Value = 500
$ws.Range("L1041").Value = 600
$ws.Range("M1041").Value = 550
$ws.Range("N1041").Value = "$/unidad"
$ws.Range("O1041").Value = "Provincia del Elquí"
$ws.Range("P1041").Value = 550
$ws.Range("Q1041").Value = 1
$ws.Range("R1041").Value = "Hortaliza"
